# The "Notes" slide (bug list: "Bugs found:" / "Bug found in APK release...")
# is no longer relevant because the underlying bugs were fixed, so the whole
# slide (and its notes page) is removed from the deck. It is slide #7 in the
# deck (between the "Architecture..." slide and the "Work done" slide).
$p = $ppt.ActivePresentation

$p.Slides.Item(7).Delete()
